$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: doctyp_code PSP -> DOC001 (rest of row unchanged)
$ws.Range("A5").Value = "DOC001"

# New rows 8-36: doctyp_code / doccat_code pairs.
# All new rows share: lang_code=ara, is_active=TRUE, cr_by=superadmin, cr_dtimes=now()
$newRows = @(
    @(8,  "DOC001", "POI"),
    @(9,  "DOC002", "POI"),
    @(10, "DOC003", "POI"),
    @(11, "DOC004", "POI"),
    @(12, "DOC005", "POI"),
    @(13, "DOC006", "POI"),
    @(14, "DOC007", "POI"),
    @(15, "DOC008", "POI"),
    @(16, "DOC009", "POI"),
    @(17, "DOC010", "POI"),
    @(18, "DOC011", "POI"),
    @(19, "DOC012", "POI"),
    @(20, "DOC001", "POA"),
    @(21, "DOC013", "POA"),
    @(22, "DOC014", "POA"),
    @(23, "DOC015", "POA"),
    @(24, "DOC004", "POA"),
    @(25, "DOC005", "POA"),
    @(26, "DOC006", "POA"),
    @(27, "DOC016", "POA"),
    @(28, "DOC017", "POA"),
    @(29, "DOC018", "POA"),
    @(30, "DOC008", "POA"),
    @(31, "DOC024", "POR"),
    @(32, "DOC025", "POR"),
    @(33, "DOC026", "POR"),
    @(34, "DOC001", "POR"),
    @(35, "DOC027", "POR"),
    @(36, "DOC028", "POR")
)

foreach ($r in $newRows) {
    $rowIdx = $r[0]
    $doctyp = $r[1]
    $doccat = $r[2]
    $ws.Cells.Item($rowIdx, 1).Value = $doctyp
    $ws.Cells.Item($rowIdx, 2).Value = $doccat
    $ws.Cells.Item($rowIdx, 3).Value = "ara"
    $ws.Cells.Item($rowIdx, 4).Value = $true
    $ws.Cells.Item($rowIdx, 5).Value = "superadmin"
    $ws.Cells.Item($rowIdx, 6).Value = "now()"
}

# Update selection to match post-edit state (selection moves past the data, to G1)
$ws.Range("G1:XFD1048576").Select()
